# "correct some spelling mistakes (from meeting with Tara)"
#
# The API-call names in column B (and a couple of matching descriptions in
# column C) contained typos ("realtionships", "wihtid", ...). This fixes
# the spelling. Excel will automatically fold the corrected strings into
# the shared-string table (reusing/dropping old entries as needed), so we
# just need to set the corrected literal text on each affected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value  = "gettradition/withid/{tradId}"
$ws.Range("B8").Value  = "deletetradition/withid/{tradId}"
$ws.Range("B9").Value  = "newtraditionwithgraphml"
$ws.Range("B10").Value = "getdot/fromtradition/{tradId}"

$ws.Range("B11").Value = "/reading"

$ws.Range("B12").Value = "getreading/fromtradition/{tradId}/withreadingid/{readId}"
$ws.Range("B13").Value = "duplicatereading/fromtradition/{tradId}"
$ws.Range("B14").Value = "mergereadings/fromtradition/{tradId}/firstReading/{firstReadId}/secondReading/{secondReadId}"
$ws.Range("B15").Value = "splitreading/fromtradition/{tradId}/ofreading/{readId}"
$ws.Range("B16").Value = "getnextreading/fromwitness/{textId}/ofreading/{readId}"
$ws.Range("B17").Value = "getpreviousreading/fromwitness/{textId}/ofreading/{readId}"
$ws.Range("B18").Value = "getallreadings/fromtradition/{tradId}"
$ws.Range("B19").Value = "getidenticalreadings/fromtradition/{tradId}/fromstartrank/{startRank}/toendrank/{endRank}"
$ws.Range("B20").Value = "couldbeidenticalreadings/fromtradition/{tradId}/fromstartrank/{startRank}/toendrank/{endRank}"
$ws.Range("B21").Value = "compressreadings/fromtradition/{tradId}/readingone/{readId1}/readingtwo/{readId2}"

$ws.Range("B23").Value = "gettext/fromtradition/{tradId}/ofwitness/{textId}"
$ws.Range("C23").Value = "getWitnessAsText"
$ws.Range("C24").Value = "getWitnessAsTextBetweenRanks"

$ws.Range("B28").Value = "createrelationship/intradition/{texId}"
$ws.Range("B29").Value = "getallrelationships/formtradition/{textId}"
$ws.Range("B30").Value = "deleterelationship/fromtradition/{textId}"
$ws.Range("B31").Value = "deleterelationshipsbyid/fromtradition/{textId}/withrealtionship/{relationshipId}"

$ws.Range("B34").Value = "createuser"
$ws.Range("B37").Value = "gettraditions/ofuser/{userId}"

$ws.Range("B39").Value = "getallstemmata/fromtradition/{tradId}"
$ws.Range("C39").Value = "getAllStemma"
$ws.Range("B40").Value = "getstemma/fromtradition/{tradId}/withtitle/{stemmaTitle}"
$ws.Range("B41").Value = "newstemma/intradition/{tradId}"
$ws.Range("B42").Value = "reorientstemma/fromtradition/{tradId}/withtitle/{stemmaTitle}/withnewrootnode/{nodeId}"
$ws.Range("C42").Value = "reorientStemma"

# These three corrections introduce brand-new shared-string entries (the
# old misspelled strings become unused and are dropped); set them last, in
# this order, so the new entries land at the tail of the shared-string
# table in the same order as upstream.
$ws.Range("B36").Value = "deleteuser/withid/{userId}"
$ws.Range("B35").Value = "getuser/withid/{userId}"
$ws.Range("B6").Value  = "getallrelationships/{tradId}"

# Cosmetic: the author's view had scrolled down a bit and moved the active
# cell/selection to B6 by the time they saved.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B6").Select()
